# Add a Bootstrap "hero" markup block (plus a trailing <script> include) as
# plain paragraphs of text, right after the existing "Our thoughtfully
# prepared..." paragraph and before the blank paragraph that already sits
# at the end of the document body (ahead of the section break).

$d = $word.ActiveDocument

# Locate the paragraph holding the tail end of the existing body copy so we
# don't have to hard-code a paragraph index.
$finder = $d.Content
$finder.Find.Execute("Our thoughtfully prepared party bags", $true, $false,
                      $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$finder.Expand(4) | Out-Null   # wdParagraph -> grow the found range to the full paragraph
$anchorStart = $finder.Start

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $anchorStart) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    # Fallback: paragraph just before the trailing blank paragraph.
    $anchorIndex = $d.Paragraphs.Count - 1
}

# The new content, one entry per paragraph to insert ("" => blank paragraph).
$newParagraphs = @(
    '',
    '',
    '<!-- Hero Section -->',
    '    <section class="hero text-center py-5">',
    '      <div class="container">',
    '        <h1 class="hero-title h2 mb-3">Welcome to Cake, Fun and Confetti</h1>',
    '        <p class="hero-subtitle lead mb-4">',
    '          We''ll help you plan and organise a perfect party, no matter how much, or how little help you need!',
    '        </p>',
    '        <a href="#contact" class="btn btn-primary">Plan your party!</a>',
    '      </div>',
    '    </section>',
    '',
    '    <!-- Bootstrap JS -->',
    '    <script src="https://cdn.jsdelivr.net/npm/bootstrap@5.3.3/dist/js/bootstrap.bundle.min.js"></script>',
    ''
)

$insertIndex = $anchorIndex
foreach ($line in $newParagraphs) {
    $r = $d.Paragraphs.Item($insertIndex).Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $insertIndex = $insertIndex + 1

    if ($line -ne "") {
        $newPara = $d.Paragraphs.Item($insertIndex)
        $newPara.Range.InsertAfter($line)
    }
}
